# Daily Report update: 2026-01-14 commit
# Adds a new date block (serial 46035 / 2026-01-13) of 22 rows to Daily_Data,
# and rolls the LOOMIS INTERNATIONAL (US) LLC Eligible RECEIVED amount
# (10177.043) through the Today_Summary and Monthly_Stats rollup sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Daily_Data: append the new date block (rows 156-177)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Daily_Data")

$newDateSerial = 46035

$rows = @(
    @("ASAHI DEPOSITORY LLC Registered", 0, 0, 0, 0, 0, 0),
    @("ASAHI DEPOSITORY LLC Eligible", 0, 0, 0, 0, 0, 0),
    @("BRINK'S, INC. Registered", 90027.72500000001, 0, 0, 0, 0, 90027.72500000001),
    @("BRINK'S, INC. Eligible", 5075.067, 0, 0, 0, 0, 5075.067),
    @("CNT DEPOSITORY, INC. Registered", 1246.06, 0, 0, 0, 0, 1246.06),
    @("CNT DEPOSITORY, INC. Eligible", 0, 0, 0, 0, 0, 0),
    @("DELAWARE DEPOSITORY Registered", 1633.941, 0, 0, 0, 0, 1633.941),
    @("DELAWARE DEPOSITORY Eligible", 18509.729, 0, 0, 0, 0, 18509.729),
    @("HSBC BANK, USA Registered", 1295.223, 0, 0, 0, 0, 1295.223),
    @("HSBC BANK, USA Eligible", 9281.978999999999, 0, 0, 0, 0, 9281.978999999999),
    @("INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered", 2395.448, 0, 0, 0, 0, 2395.448),
    @("INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible", 0, 0, 0, 0, 0, 0),
    @("JP MORGAN CHASE BANK NA Registered", 124991.729, 0, 0, 0, 0, 124991.729),
    @("JP MORGAN CHASE BANK NA Eligible", 125407.673, 0, 0, 0, 0, 125407.673),
    @("LOOMIS INTERNATIONAL (US) LLC Registered", 68084.33, 0, 0, 0, 0, 68084.33),
    @("LOOMIS INTERNATIONAL (US) LLC Eligible", 106188.481, 10177.043, 0, 10177.043, 0, 116365.524),
    @("MALCA-AMIT USA, LLC Registered", 395.145, 0, 0, 0, 0, 395.145),
    @("MALCA-AMIT USA, LLC Eligible", 0, 0, 0, 0, 0, 0),
    @("MANFRA, TORDELLA & BROOKES, LLC Registered", 54605.27, 0, 0, 0, 0, 54605.27),
    @("MANFRA, TORDELLA & BROOKES, LLC Eligible", 1068.408, 0, 0, 0, 0, 1068.408),
    @("STONEX PRECIOUS METALS LLC Registered", 14122.765, 0, 0, 0, 0, 14122.765),
    @("STONEX PRECIOUS METALS LLC Eligible", 16.075, 0, 0, 0, 0, 16.075)
)

$r = 156
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 1).Value = $newDateSerial
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. Today_Summary: roll the new RECEIVED amount into LOOMIS totals (row 9)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Today_Summary")
$ws2.Cells.Item(9, 2).Value = 116365.524   # Eligible
$ws2.Cells.Item(9, 4).Value = 184449.854   # Total_Stock

# ---------------------------------------------------------------------------
# 3. Monthly_Stats: roll the new RECEIVED amount into the month's grand
#    totals (row 2) and into the LOOMIS Eligible detail row (row 21)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Monthly_Stats")
$ws3.Cells.Item(2, 2).Value = 275724.455   # Eligible
$ws3.Cells.Item(2, 4).Value = 634522.091   # Grand_Total

$ws3.Cells.Item(21, 3).Value = 10177.043    # RECEIVED
$ws3.Cells.Item(21, 5).Value = 116365.524   # TOTAL_TODAY
